# Applies a cyclic rotation of the species-observation data held in rows 2-4:
#   old row 4 -> new row 2
#   old row 2 -> new row 3
#   old row 3 -> new row 4
# Only the columns that actually differ between the three rows are touched
# (A, B, D, E, F, G, H, Q, R, Z, AB, AC); columns that are identical across
# the three rows (C, I, K, P, S, T, U, V, W, Y, AD, AE, AG, AT, AW, AX, AY)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB", "AC")

# Snapshot the current ("before") values of rows 2, 3 and 4 for the columns
# that participate in the rotation.
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("$col" + "2").Value()
    $row3[$col] = $ws.Range("$col" + "3").Value()
    $row4[$col] = $ws.Range("$col" + "4").Value()
}

# Write the rotated values: row4 -> row2, row2 -> row3, row3 -> row4.
foreach ($col in $cols) {
    $ws.Range("$col" + "2").Value = $row4[$col]
    $ws.Range("$col" + "3").Value = $row2[$col]
    $ws.Range("$col" + "4").Value = $row3[$col]
}

# Row 2 gains the public-comment text previously on row 4, and rows 3/4 no
# longer carry a comment - clear the now-stale AC value on row 4 explicitly
# (it was already copied above, but make sure row 4 ends up blank since the
# source row 3 never had one).
if ([string]::IsNullOrEmpty($row3["AC"])) {
    $ws.Range("AC4").ClearContents()
}
